$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2..221)
# from 45184 to 45186 (serial date values).
for ($r = 2; $r -le 221; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45186
    }
}

# Add a friendly display name ("A 45791-2019") as the second argument of the
# HYPERLINK formulas in row 2 (columns S, T, V, W, X, Y).
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HERRLJUNGA/artfynd/A 45791-2019.xlsx", "A 45791-2019")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HERRLJUNGA/kartor/A 45791-2019.png", "A 45791-2019")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HERRLJUNGA/klagomål/A 45791-2019.docx", "A 45791-2019")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HERRLJUNGA/klagomålsmail/A 45791-2019.docx", "A 45791-2019")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HERRLJUNGA/tillsyn/A 45791-2019.docx", "A 45791-2019")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HERRLJUNGA/tillsynsmail/A 45791-2019.docx", "A 45791-2019")'

$wb.Save()
